# Update calibration parameters and abs_epsi_autocorr series for each model sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.01109759028167569
$ws.Range("B3").Value = 0.7954978562199575
$ws.Range("B4").Value = 3.467916129465851
$ws.Range("B5").Value = '[1.0, 0.3025720289175844, 0.2721670726266026, 0.25387086147244425, 0.30050729359130185, 0.4410391118631507, 0.26534041264157804, 0.2732724434404467, 0.22440448698739146, 0.24557266031483727, 0.206327595248993, 0.2295578942951122, 0.21843307567942796, 0.2319503637838938, 0.2222817049018749, 0.242291518635085, 0.21447296949601752, 0.22601426320293905, 0.20917647153927402, 0.19866297734490548]'

$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 0.9142000749344323
$ws.Range("B4").Value = -0.5499784532826222
$ws.Range("B5").Value = 0.7409012743822112
$ws.Range("B6").Value = 2.702539947776487
$ws.Range("B7").Value = 0.5588548102736736
$ws.Range("B8").Value = 0.7266939026303234
$ws.Range("B9").Value = 2.035209533783421
$ws.Range("B10").Value = '[1.0, 0.30828147383071797, 0.2740333890296752, 0.2557075253571407, 0.28732205052130044, 0.310155596599058, 0.23979776807310726, 0.24939082546093957, 0.21098757855133077, 0.1897614264675183, 0.17233899140388115, 0.20862399803657924, 0.18794873846187962, 0.23474574598500886, 0.1835548567183007, 0.21264135532071696, 0.19388997426013266, 0.19400845199018657, 0.1739831612513745, 0.1772967009857266]'

$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.004917013543116382
$ws.Range("B3").Value = 0.01382885587038979
$ws.Range("B4").Value = 0.08404298000164705
$ws.Range("B5").Value = 0.9159570199779029
$ws.Range("B6").Value = '[0.9999999999999999, 0.014215798214269885, -0.032768154212978025, -0.04873804224200899, -0.02553423179016087, 0.27592921156337985, 0.0169560882247787, 0.01944603049899669, -0.017690299993236593, -0.016038632757813334, -0.07182011410421124, -0.04675532592617912, -0.030468609700310293, -0.022758209798758023, 0.03688400244989628, 0.004667543669559729, -0.005554312904124156, -0.057569123755357235, -0.0378993359352002, -0.02451952382228693]'

$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.01148342628388638
$ws.Range("B3").Value = 0.01394289010722019
$ws.Range("B4").Value = 0.07129698137505622
$ws.Range("B5").Value = 0.9173905446124516
$ws.Range("B6").Value = '[0.9999999999999999, 0.01498808057504423, -0.032586889522347874, -0.047963438198307715, -0.024673982274183617, 0.27477191812059043, 0.017410278892811937, 0.019537577642349934, -0.017604227812687466, -0.015559555206280249, -0.07115955449377741, -0.046284222786104424, -0.029810744557196536, -0.022429057744130152, 0.03759326794842329, 0.005345840748064989, -0.005354193194266831, -0.05806082329755134, -0.0372930962233834, -0.024651293113016035]'
$ws.Range("B7").Value = 0.02262496324655447

$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.001021012211469371
$ws.Range("B3").Value = 0.009982672908142485
$ws.Range("B4").Value = 0.09398901493423498
$ws.Range("B5").Value = 0.9187737720627062
$ws.Range("B6").Value = '[1.0, 0.03225010707379681, -0.008094478047663346, -0.01950175328116482, 0.01793112570211922, 0.20820212472771965, 0.012480849181006164, 0.011057803170610744, -0.03046301382105047, -0.022026219216421592, -0.055040302910202045, -0.041447480072262684, -0.04218421371892634, -0.035872962545917576, 0.02485395409233621, -0.00489371330339512, -0.014304444734273407, -0.03961969261708188, -0.029916716436765887, -0.029049664088345084]'
$ws.Range("B7").Value = -0.02552556613595806
$ws.Range("B9").Value = 0.7838241161838281
